{"js": "// Datenblatt verbessert und ans Frontend angebunden\n//\n// 1) The \"Aktueller Wert: ${Aktueller Wert}\" paragraph gets explicit\n//    paragraph spacing (before = 0pt, after = 7pt == 140 twips).\n// 2) Every paragraph that followed it (the empty spacer paragraph, the\n//    ${block_Freitexte} / ${Titel} / ${Beschreibung} paragraphs, another\n//    empty spacer, and ${/block_Freitexte}) is removed from the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Aktueller Wert:\";\nlet markerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error(`Could not find paragraph containing \"${marker}\"`);\n}\n\n// Give the \"Aktueller Wert\" paragraph its new spacing (values are in\n// points; 0pt before / 7pt after == w:spacing w:before=\"0\" w:after=\"140\").\nconst targetParagraph = paragraphs.items[markerIndex];\ntargetParagraph.spaceBefore = 0;\ntargetParagraph.spaceAfter = 7;\n\n// Delete every paragraph that comes after it (the old ${block_Freitexte}\n// block and its surrounding blank paragraphs).\nfor (let i = paragraphs.items.length - 1; i > markerIndex; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Datenblatt verbessert und ans Frontend angebunden\n#\n# 1) The \"Aktueller Wert: ${Aktueller Wert}\" paragraph gets explicit\n#    paragraph spacing (before = 0pt, after = 7pt == 140 twips).\n# 2) Every paragraph that followed it (the empty spacer paragraph, the\n#    ${block_Freitexte} / ${Titel} / ${Beschreibung} paragraphs, another\n#    empty spacer, and ${/block_Freitexte}) is removed from the document.\n\n$d = $word.ActiveDocument\n\n$marker = \"Aktueller Wert:\"\n$markerIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*$marker*\") {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -eq -1) {\n    throw \"Could not find paragraph containing '$marker'\"\n}\n\n# Give the \"Aktueller Wert\" paragraph its new spacing (values are in\n# points; 0pt before / 7pt after == w:spacing w:before=\"0\" w:after=\"140\").\n$target = $d.Paragraphs.Item($markerIndex)\n$target.Range.ParagraphFormat.SpaceBefore = 0\n$target.Range.ParagraphFormat.SpaceAfter = 7\n\n# Delete every paragraph that comes after it (the old ${block_Freitexte}\n# block and its surrounding blank paragraphs), from the end backwards so\n# indices stay valid.\nfor ($i = $d.Paragraphs.Count; $i -gt $markerIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
